$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.511.21"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.77%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.90"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.12%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.03%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.38"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.57%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6505"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +2.08%  "

$ws.Range("E7").Value = "  +0.09%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.89"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +4.10%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07502"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -0.33%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2982"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +0.19%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.46"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +1.18%  "

$ws.Range("E12").Value = "  -0.64%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.857.34"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +0.54%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.029"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.56%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6856"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.33%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.62"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.07%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009468"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.40%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.119"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.63%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.565.46"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.45%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.105.39"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.60%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "237.90"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.68%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.61"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.04%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.701"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +4.69%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.001"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.03%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.43"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -1.39%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1418"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.23%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.513"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.26%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.84"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -0.63%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06074"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +0.94%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.490"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.73%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.261"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -0.20%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.142"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +0.06%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.070"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -1.84%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.188"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +3.29%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.860"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -1.03%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7258"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.01%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.595"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.55%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.794"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -2.34%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01782"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +0.12%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.201.44"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.08%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.228"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -2.04%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9070"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("E44").Value = "  -0.04%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.016.10"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +0.37%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.72"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -0.53%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.55"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +0.16%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.443"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +10.90%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000123"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +0.00%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4052"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -0.79%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.152"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -1.98%  "
